# Update the "cryptos" price/volume snapshot with freshly scraped values.
# Mirrors the GitHub Actions job that refreshes this sheet on a schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Text
    )

    $rng = $ws.Range($Cell)

    # Some of the "Price" values (column D) are plain decimal numbers
    # (e.g. "244.39"). If we assign them to .Value directly, Excel's COM
    # layer auto-coerces the string into a floating point number, which
    # both loses the exact textual representation and introduces binary
    # floating point noise (244.39 -> 244.38999999999999...). Forcing the
    # cell to Text format for the duration of the assignment keeps the
    # value as the literal string we want; ClearFormats() afterwards
    # removes the temporary formatting again so the cell is left without
    # any explicit style, matching its original (unstyled) state.
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.ClearFormats()
}

# --- Bitcoin ---
Set-TextValue "D2" "29.890.59"

# --- Ethereum ---
Set-TextValue "D3" "1.894.81"

# --- TetherUSD ---
$ws.Range("E4").Value = "  -0.04%  "

# --- XRP ---
$ws.Range("E5").Value = "  -2.26%  "

# --- BNB ---
Set-TextValue "D6" "244.39"
$ws.Range("E6").Value = "  +0.36%  "

# --- USDC ---
$ws.Range("E7").Value = "  -0.06%  "

# --- Cardano ---
Set-TextValue "D8" "0.3124"
$ws.Range("E8").Value = "  -0.74%  "

# --- Solana ---
$ws.Range("E9").Value = "  +1.27%  "

# --- Dogecoin ---
Set-TextValue "D10" "0.07224"
$ws.Range("E10").Value = "  -0.18%  "

# --- TRON ---
Set-TextValue "D11" "0.08896"
$ws.Range("E11").Value = "  +9.77%  "

# --- Polygon ---
Set-TextValue "D12" "0.7721"
$ws.Range("E12").Value = "  +0.90%  "

# --- Polkadot ---
Set-TextValue "D13" "5.428"
$ws.Range("E13").Value = "  -2.68%  "

# --- WrappedEther ---
Set-TextValue "D14" "1.879.57"
$ws.Range("E14").Value = "  -0.35%  "

# --- Uniswap ---
$ws.Range("E16").Value = "  +0.09%  "

# --- WrappedBTC ---
Set-TextValue "D17" "29.861.00"
$ws.Range("E17").Value = "  -0.06%  "

# --- Avalanche ---
$ws.Range("E18").Value = "  +0.14%  "

# --- BitcoinCash ---
Set-TextValue "D19" "245.73"
$ws.Range("E19").Value = "  +0.72%  "

# --- ShibaInu ---
Set-TextValue "D20" "0.000007853"
$ws.Range("E20").Value = "  +0.81%  "

# --- Dai ---
Set-TextValue "D21" "0.9998"
$ws.Range("E21").Value = "  -0.18%  "

# --- Chainlink ---
Set-TextValue "D22" "8.124"
$ws.Range("E22").Value = "  -0.83%  "

# --- WrappedliquidstakedEther2.0 ---
Set-TextValue "D23" "2.128.02"
$ws.Range("E23").Value = "  -2.22%  "

# --- BinanceUSD ---
$ws.Range("E24").Value = "  -0.04%  "

# --- Stellar ---
Set-TextValue "D25" "0.1588"
$ws.Range("E25").Value = "  -4.06%  "

# --- Cosmos ---
Set-TextValue "D26" "9.513"
$ws.Range("E26").Value = "  +1.11%  "

# --- Monero ---
Set-TextValue "D27" "162.55"
$ws.Range("E27").Value = "  -0.74%  "

# --- LidoDAOToken ---
Set-TextValue "D29" "2.042"
$ws.Range("E29").Value = "  -0.93%  "

# --- Toncoin ---
$ws.Range("E30").Value = "  +1.94%  "

# --- PancakeSwap ---
Set-TextValue "D31" "1.543"
$ws.Range("E31").Value = "  -0.13%  "

# --- Filecoin ---
Set-TextValue "D32" "4.557"
$ws.Range("E32").Value = "  +1.88%  "

# --- InternetComputer(DFINITY) ---
$ws.Range("E33").Value = "  +0.31%  "

# --- Hedera ---
Set-TextValue "D34" "0.05498"
$ws.Range("E34").Value = "  -0.70%  "

# --- ARBITRUM ---
Set-TextValue "D35" "1.249"
$ws.Range("E35").Value = "  -1.76%  "

# --- ImmutableX ---
Set-TextValue "D36" "0.7492"
$ws.Range("E36").Value = "  +1.36%  "

# --- Frax ---
Set-TextValue "D37" "0.9982"
$ws.Range("E37").Value = "  -0.25%  "

# --- HuobiToken ---
Set-TextValue "D38" "2.711"
$ws.Range("E38").Value = "  +3.49%  "

# --- VeChain ---
Set-TextValue "D39" "0.01955"
$ws.Range("E39").Value = "  +1.60%  "

# --- MXToken ---
Set-TextValue "D40" "2.786"
$ws.Range("E40").Value = "  +0.25%  "

# --- TheSandbox ---
Set-TextValue "D41" "0.4500"
$ws.Range("E41").Value = "  +1.78%  "

# --- Aave ---
Set-TextValue "D42" "73.77"
$ws.Range("E42").Value = "  -0.51%  "

# --- Maker ---
Set-TextValue "D43" "1.090.31"
$ws.Range("E43").Value = "  -4.76%  "

# --- FraxShare ---
Set-TextValue "D44" "6.029"
$ws.Range("E44").Value = "  +2.43%  "

# --- TrustWalletToken ---
Set-TextValue "D45" "0.8549"
$ws.Range("E45").Value = "  +0.60%  "

# --- PaxDollar ---
$ws.Range("E46").Value = "  -0.05%  "

# --- RenderToken ---
Set-TextValue "D47" "1.883"
$ws.Range("E47").Value = "  +0.28%  "

# --- Quant ---
Set-TextValue "D48" "102.48"
$ws.Range("E48").Value = "  -1.73%  "

# --- Aptos ---
Set-TextValue "D49" "7.613"
$ws.Range("E49").Value = "  +2.18%  "

# --- EnergySwap ---
Set-TextValue "D50" "9.885"
$ws.Range("E50").Value = "  -1.19%  "

# --- SynthetixNetwork ---
Set-TextValue "D51" "2.967"
$ws.Range("E51").Value = "  -1.43%  "
